# Re-order the header columns in row 2 of the Staging.MilestoneValues
# template: move the ID column and all *BusinessKey columns to the front
# (alphabetically amongst themselves), followed by the remaining columns
# (also alphabetically amongst themselves) - matching the "re-ordered
# columns" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newOrder = @(
    "MilestoneValues_ID",
    "AgeBandBusinessKey",
    "BusinessKey",
    "CommunityTypeBusinessKey",
    "DataVersionBusinessKey",
    "DonorBusinessKey",
    "FrameworkBusinessKey",
    "GenderBusinessKey",
    "GroupBusinessKey",
    "InstitutionBusinessKey",
    "LocationBusinessKey",
    "MilestoneBusinessKey",
    "OrganizationBusinessKey",
    "ReportingPeriodBusinessKey",
    "ResultAreaBusinessKey",
    "StrategicElementBusinessKey",
    "ActualDate",
    "ActualLabel",
    "ActualValue",
    "GroupVersion",
    "Notes"
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $newOrder[$i]
}
